$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old block (rows 32:44) entirely so no stale row formatting remains ---
$ws.Range("A32:A44").EntireRow.Delete()

# --- Row 32 ---
$ws.Range("A32").Value = 'Save system base types when the project is saved, updating their former record so that their id stays the same. But deleting and then inserting of methods, props and events is the way to go for them.'
$ws.Range("A32").Font.Strikethrough = $true
$ws.Range("A32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 58

# --- Row 33 ---
$ws.Range("A33").Value = 'Don''t forget to delete m, p, e from updated SBTs.'
$ws.Range("A33").Font.Strikethrough = $true
$ws.Range("A33").WrapText = $true

# --- Row 34 ---
$ws.Range("A34").Value = 'Why did it work?'
$ws.Range("A34").Font.Strikethrough = $true
$ws.Range("A34").WrapText = $true
$ws.Range("B34").Value = 'I should have added new Types twice.'
$ws.Range("B34").Font.Strikethrough = $true
$ws.Range("B34").WrapText = $true

# --- Row 35 ---
$ws.Range("A35").Value = 'Any type that''s a base type (even one that''s not new--if not an SBT) will get a new id when it''s written in save or save as and, thus the guys who point to him has to have baseTypeId updated.'
$ws.Range("A35").Font.Strikethrough = $true
$ws.Range("A35").WrapText = $true
$ws.Range("B35").Value = 'I''m going to need to mark types that are base types (not SBTs) si that I can update their derived types'' bastTypeIds.'
$ws.Range("B35").Font.Strikethrough = $true
$ws.Range("B35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 58

# --- Row 36 ---
$ws.Range("A36").Font.Strikethrough = $true
$ws.Range("A36").WrapText = $true
$ws.Range("B36").Font.Strikethrough = $true
$ws.Range("B36").WrapText = $true

# --- Row 37 ---
$ws.Range("A37").Font.Strikethrough = $true
$ws.Range("A37").WrapText = $true
$ws.Range("B37").Font.Strikethrough = $true
$ws.Range("B37").WrapText = $true

# --- Row 38 ---
$ws.Range("A38").Font.Strikethrough = $true
$ws.Range("A38").WrapText = $true

# --- Row 39 ---
$ws.Range("A39").Value = 'TO DO'
$ws.Range("A39").Font.Bold = $true
$ws.Range("A39").Font.Size = 14
$ws.Range("A39").HorizontalAlignment = -4108
$ws.Range("A39").WrapText = $true
$ws.Range("B39").Font.Strikethrough = $true
$ws.Range("B39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 18.5

# --- Row 40 ---
$ws.Range("A40").Value = 'Rename TI in PropertyGrid'
$ws.Range("A40").WrapText = $true
$ws.Range("B40").Value = 'Still broken?'
$ws.Range("B40").WrapText = $true

# --- Row 41 ---
$ws.Range("A41").Value = 'All Projects menu items are available after closing a Project. Also, closing a project has to clear the browser tab.'
$ws.Range("A41").WrapText = $true
$ws.Range("B41").Value = 'This is caused by a complex chain of events that I have to step through.'
$ws.Range("B41").WrapText = $true
$ws.Rows.Item(41).RowHeight = 43.5

# --- Row 42 ---
$ws.Range("A42").Value = 'Reverse enabled/disabled menu colors.'
$ws.Range("A42").WrapText = $true

# --- Row 43 ---
$ws.Range("A43").Value = 'See if I can auto-update the system base type SQL script when a project with system base types maintenance is saved.'
$ws.Range("A43").WrapText = $true
$ws.Range("B43").Value = 'Or always, if I don''t know.'
$ws.Range("B43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 43.5

# --- Row 44 ---
$ws.Range("A44").Value = 'In TypeWell: Delete current type should be disabled for: App Type; any SBT; any Type in the current Comic that is a base type for another type in that comic.'
$ws.Range("A44").WrapText = $true
$ws.Range("B44").Value = 'Think more about deleting SBTs beyond the first 5.'
$ws.Range("B44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 43.5

# --- Row 45 ---
$ws.Range("A45").Value = 'A New SBT should probably require an image.'
$ws.Range("A45").WrapText = $true
$ws.Range("B45").Value = 'Discuss w/Ken. He agrees, but, since it''s only us who add them, says the program doesn''t have to enforce. Maybe.'
$ws.Range("B45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 29

# --- Row 46 ---
$ws.Range("A46").Value = 'Add  new fields to ProjectBO.routeRetrieveType and its 3 arrays. And to routeRetrieveMethod.'
$ws.Range("A46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 29

# --- Row 47 ---
$ws.Range("A47").Value = 'If a non-SBT base type is deleted, either the type(s) that point to it have to be nulled, or better, it should be prohibited in the UI.'
$ws.Range("A47").WrapText = $true
$ws.Rows.Item(47).RowHeight = 43.5

# --- Row 48 ---
$ws.Range("A48").Value = 'Add usergroups'
$ws.Range("A48").WrapText = $true

# --- Selection / view ---
$ws.Range("A47").Select()

